$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: sort-order id bumped by one ---
$ws.Range("B11").Value = 93134

# --- Rows 13 & 14: the two observation records swap places ---
# Use a scratch row far outside the used range to stage row 13's original
# content while row 14's content is copied into row 13, then move the
# staged content into row 14. Using whole-row copy/paste (instead of
# cell-by-cell value assignment) preserves each cell's original storage
# type (numeric vs text) and avoids introducing new number formats/styles
# for text-like numeric strings (e.g. "1", "5" in column I).
$scratchRow = 500

$ws.Rows.Item(13).Copy() | Out-Null
$ws.Rows.Item($scratchRow).PasteSpecial(-4104) | Out-Null   # xlPasteAll

$ws.Rows.Item(14).Copy() | Out-Null
$ws.Rows.Item(13).PasteSpecial(-4104) | Out-Null

$ws.Rows.Item($scratchRow).Copy() | Out-Null
$ws.Rows.Item(14).PasteSpecial(-4104) | Out-Null

$ws.Rows.Item($scratchRow).Clear() | Out-Null

# The whole-row swap leaves stray cells behind in columns that aren't
# shared between the two records' column layouts (row 13 used J, row 14
# used K/M) because paste doesn't blank out destination cells beyond
# what changed. Clean those up explicitly.
$ws.Range("J13").Value = $null
$ws.Range("K14").Value = $null
$ws.Range("M14").Value = $null

# Column B ("Taxonsorteringsordning") doesn't follow the row swap like
# the rest of the record - the new row 13 keeps the sort id that used to
# belong to row 14 (already true after the row swap above), but the new
# row 14 gets a freshly minted id rather than reusing row 13's old one.
$ws.Range("B14").Value = 93096

Write-Output "done"
